# Fix alternating column formatting (with the "buffer"/off-by-one bug on the
# last column) on the Count sheet, and correct the report date in both sheet
# names from 09-02-2022 (MM-DD) to 02-09-2022 (DD-MM).

$wb = $excel.ActiveWorkbook

$wsData  = $wb.Worksheets.Item(1)
$wsCount = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------
# 1. Rename the sheets
# ---------------------------------------------------------------------
$wsData.Name  = "MaddenCo Data 02-09-2022"
$wsCount.Name = "MaddenCo Count 02-09-2022"

# ---------------------------------------------------------------------
# 2. Re-color the alternating columns on the Count sheet.
#    Previously columns alternated B/D/F/H vs C/E/G, but column H ended up
#    sharing the "odd" color with G instead of matching B/D/F - this is the
#    formatting bug being fixed here. We also refresh the palette.
# ---------------------------------------------------------------------

# Template cell that already has the right font/border/alignment
# (bold-free font, thin border, centered) used for the data rows.
$dataTemplate = $wsCount.Range("C2")

# New "odd" data columns (C, E, G) -> light blue DAEEF3
foreach ($col in @("C", "E", "G")) {
    $target = $wsCount.Range("$col" + "2:" + "$col" + "5")
    $dataTemplate.Copy()
    $target.PasteSpecial(-4122)
    $target.Interior.Color = 0xF3EEDA
}

# New "even" data columns (B, D, F, H) -> light peach FDE9D9
# (H now correctly matches B/D/F, fixing the old buffer bug)
foreach ($col in @("B", "D", "F", "H")) {
    $target = $wsCount.Range("$col" + "2:" + "$col" + "5")
    $dataTemplate.Copy()
    $target.PasteSpecial(-4122)
    $target.Interior.Color = 0xD9E9FD
}

# Template cell for the header row (row 1) that already has the right
# font/border/alignment.
$headerTemplate = $wsCount.Range("B1")

# Header "even" columns (B, D, F, H) -> orange F79646
foreach ($col in @("B", "D", "F", "H")) {
    $target = $wsCount.Range("$col" + "1")
    $headerTemplate.Copy()
    $target.PasteSpecial(-4122)
    $target.Interior.Color = 0x4696F7
}

# Header "odd" columns (C, E, G) -> teal 4BACC6
foreach ($col in @("C", "E", "G")) {
    $target = $wsCount.Range("$col" + "1")
    $headerTemplate.Copy()
    $target.PasteSpecial(-4122)
    $target.Interior.Color = 0xC6AC4B
}

$excel.CutCopyMode = $false
